# Auto-generated Excel COM-interop script applying the Famfrit_Profits price/profit updates.
# For each affected row (keyed by the Leve Item ID in column G) across the ALC/ARM/BSM/CRP/
# CUL/GSM/LTW/WVR sheets, refresh the market-price-derived columns (H:N).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 500
$ws.Range("J49").Value = 500
$ws.Range("L49").Value = 1500
$ws.Range("N49").Value = -1772
$ws.Range("H62").Value = 1816.6666
$ws.Range("I62").Value = 975
$ws.Range("K62").Value = 975
$ws.Range("M62").Value = -351
$ws.Range("H65").Value = 1816.6666
$ws.Range("I65").Value = 975
$ws.Range("K65").Value = 4875
$ws.Range("M65").Value = -1755
$ws.Range("H125").Value = 6102.857
$ws.Range("J125").Value = 7944
$ws.Range("L125").Value = 71496
$ws.Range("N125").Value = -76416
$ws.Range("H132").Value = 6801.0527
$ws.Range("I132").Value = 7101.1113
$ws.Range("K132").Value = 21303.3339
$ws.Range("M132").Value = -18773.3339
$ws.Range("H138").Value = 5730.3145
$ws.Range("I138").Value = 1372.1818
$ws.Range("J138").Value = 7727.7915
$ws.Range("K138").Value = 4116.5454
$ws.Range("L138").Value = 23183.3745
$ws.Range("M138").Value = 1023.4546
$ws.Range("N138").Value = -33463.37450000001
$ws.Range("H139").Value = 159333.33
$ws.Range("J139").Value = 179000
$ws.Range("L139").Value = 179000
$ws.Range("N139").Value = -189280

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 33334992
$ws.Range("I74").Value = 40000972
$ws.Range("K74").Value = 40000972
$ws.Range("M74").Value = -40000098
$ws.Range("H77").Value = 33334992
$ws.Range("I77").Value = 40000972
$ws.Range("K77").Value = 200004860
$ws.Range("M77").Value = -200000492

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 807
$ws.Range("J15").Value = 807
$ws.Range("L15").Value = 807
$ws.Range("N15").Value = -1261
$ws.Range("H19").Value = 2558.9167
$ws.Range("I19").Value = 2849.7
$ws.Range("J19").Value = 1105
$ws.Range("K19").Value = 2849.7
$ws.Range("L19").Value = 1105
$ws.Range("M19").Value = -2676.7
$ws.Range("N19").Value = -1451
$ws.Range("H134").Value = 4691.2
$ws.Range("I134").Value = 4640.8335
$ws.Range("K134").Value = 13922.5005
$ws.Range("M134").Value = -11387.5005

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 332.33334
$ws.Range("I12").Value = 332.33334
$ws.Range("K12").Value = 332.33334
$ws.Range("M12").Value = -162.33334
$ws.Range("H16").Value = 5331.6665
$ws.Range("I16").Value = 5498.25
$ws.Range("K16").Value = 5498.25
$ws.Range("M16").Value = -5211.25
$ws.Range("H22").Value = 5085.1665
$ws.Range("I22").Value = 5895.8887
$ws.Range("K22").Value = 5895.8887
$ws.Range("M22").Value = -5545.8887
$ws.Range("H58").Value = 1311.25
$ws.Range("I58").Value = 1311.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1311.25
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1108.25
$ws.Range("N58").ClearContents()
$ws.Range("H86").Value = 3946.25
$ws.Range("I86").Value = 3930
$ws.Range("K86").Value = 3930
$ws.Range("M86").Value = -2807
$ws.Range("H89").Value = 3946.25
$ws.Range("I89").Value = 3930
$ws.Range("K89").Value = 19650
$ws.Range("M89").Value = -14034
$ws.Range("H99").Value = 10376.4375
$ws.Range("I99").Value = 7499.5557
$ws.Range("J99").Value = 11502.174
$ws.Range("K99").Value = 7499.5557
$ws.Range("L99").Value = 11502.174
$ws.Range("M99").Value = -6001.5557
$ws.Range("N99").Value = -14498.174
$ws.Range("H113").Value = 5331.6665
$ws.Range("I113").Value = 5498.25
$ws.Range("K113").Value = 5498.25
$ws.Range("M113").Value = -3328.25
$ws.Range("H126").Value = 10376.4375
$ws.Range("I126").Value = 7499.5557
$ws.Range("J126").Value = 11502.174
$ws.Range("K126").Value = 22498.6671
$ws.Range("L126").Value = 34506.522
$ws.Range("M126").Value = -20028.6671
$ws.Range("N126").Value = -39446.522
$ws.Range("H131").Value = 24865.666
$ws.Range("I131").Value = 13000
$ws.Range("K131").Value = 13000
$ws.Range("M131").Value = -7960
$ws.Range("H136").Value = 1311.25
$ws.Range("I136").Value = 1311.25
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3933.75
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1383.75
$ws.Range("N136").ClearContents()
$ws.Range("H141").Value = 103894.664
$ws.Range("J141").Value = 113131.5
$ws.Range("L141").Value = 113131.5
$ws.Range("N141").Value = -123491.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 40631524
$ws.Range("I4").Value = 57080150
$ws.Range("K4").Value = 171240450
$ws.Range("M4").Value = -171240338
$ws.Range("H75").Value = 6960.7144
$ws.Range("J75").Value = 1745
$ws.Range("L75").Value = 5235
$ws.Range("N75").Value = -7231
$ws.Range("H78").Value = 6960.7144
$ws.Range("J78").Value = 1745
$ws.Range("L78").Value = 15705
$ws.Range("N78").Value = -25689
$ws.Range("H86").Value = 1044.6154
$ws.Range("I86").Value = 734.44446
$ws.Range("K86").Value = 2203.33338
$ws.Range("M86").Value = -1017.33338
$ws.Range("H89").Value = 1044.6154
$ws.Range("I89").Value = 734.44446
$ws.Range("K89").Value = 6610.00014
$ws.Range("M89").Value = -682.0001400000001
$ws.Range("H97").Value = 1750.4286
$ws.Range("J97").Value = 1549.8
$ws.Range("L97").Value = 4649.4
$ws.Range("N97").Value = -5641.4
$ws.Range("H131").Value = 1318.2
$ws.Range("I131").Value = 952.25
$ws.Range("J131").Value = 1656
$ws.Range("K131").Value = 2856.75
$ws.Range("L131").Value = 4968
$ws.Range("M131").Value = 2183.25
$ws.Range("N131").Value = -15048

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7100
$ws.Range("I43").Value = 7100
$ws.Range("K43").Value = 7100
$ws.Range("M43").Value = -6949
$ws.Range("H51").Value = 59932.668
$ws.Range("I51").Value = 63199.668
$ws.Range("K51").Value = 63199.668
$ws.Range("M51").Value = -62690.668
$ws.Range("H80").Value = 7087.6665
$ws.Range("I80").Value = 6259.5835
$ws.Range("J80").Value = 7915.75
$ws.Range("K80").Value = 6259.5835
$ws.Range("L80").Value = 7915.75
$ws.Range("M80").Value = -5261.5835
$ws.Range("N80").Value = -9911.75
$ws.Range("H83").Value = 7087.6665
$ws.Range("I83").Value = 6259.5835
$ws.Range("J83").Value = 7915.75
$ws.Range("K83").Value = 31297.9175
$ws.Range("L83").Value = 39578.75
$ws.Range("M83").Value = -26305.9175
$ws.Range("N83").Value = -49562.75

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1949.8889
$ws.Range("I22").Value = 2067.6667
$ws.Range("J22").Value = 1832.1111
$ws.Range("K22").Value = 2067.6667
$ws.Range("L22").Value = 1832.1111
$ws.Range("M22").Value = -1772.6667
$ws.Range("N22").Value = -2422.1111
$ws.Range("H27").Value = 1949.8889
$ws.Range("I27").Value = 2067.6667
$ws.Range("J27").Value = 1832.1111
$ws.Range("K27").Value = 2067.6667
$ws.Range("L27").Value = 1832.1111
$ws.Range("M27").Value = -1960.6667
$ws.Range("N27").Value = -2046.1111
$ws.Range("H40").Value = 5324.95
$ws.Range("I40").Value = 5388.1763
$ws.Range("K40").Value = 5388.1763
$ws.Range("M40").Value = -5252.1763
$ws.Range("H45").Value = 25997.5
$ws.Range("I45").Value = 14500
$ws.Range("J45").Value = 37495
$ws.Range("K45").Value = 14500
$ws.Range("L45").Value = 37495
$ws.Range("M45").Value = -14093
$ws.Range("N45").Value = -38309
$ws.Range("H48").Value = 36831.668
$ws.Range("I48").Value = 36000
$ws.Range("J48").Value = 37247.5
$ws.Range("K48").Value = 36000
$ws.Range("L48").Value = 37247.5
$ws.Range("M48").Value = -35339
$ws.Range("N48").Value = -38569.5
$ws.Range("H61").Value = 2770
$ws.Range("I61").Value = 2445.077
$ws.Range("K61").Value = 2445.077
$ws.Range("M61").Value = -2243.077
$ws.Range("H82").Value = 2411.875
$ws.Range("J82").Value = 1250.75
$ws.Range("L82").Value = 1250.75
$ws.Range("N82").Value = -1972.75
$ws.Range("H85").Value = 2411.875
$ws.Range("J85").Value = 1250.75
$ws.Range("L85").Value = 1250.75
$ws.Range("N85").Value = -3746.75
$ws.Range("H113").Value = 2770
$ws.Range("I113").Value = 2445.077
$ws.Range("K113").Value = 2445.077
$ws.Range("M113").Value = -275.0770000000002
$ws.Range("H127").Value = 105000
$ws.Range("J127").Value = 105000
$ws.Range("L127").Value = 105000
$ws.Range("N127").Value = -114920

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2075.2354
$ws.Range("I107").Value = 1497
$ws.Range("J107").Value = 2199.1428
$ws.Range("K107").Value = 4491
$ws.Range("L107").Value = 6597.428400000001
$ws.Range("M107").Value = -2571
$ws.Range("N107").Value = -10437.4284

